# The author selected AF2:AI71 (the "QTD ENVIADA" / "OBS ROMANEIO" /
# "PESO Unitario" / "AREA M2 Unitario" helper columns on the "Mapa de
# Suportes" sheet) and cleared their contents, then re-uploaded the
# workbook. The previously-static 0 values go away (cells keep their
# style/number format, but lose the cached <v>0</v>), and the sheet's
# saved selection is left sitting on that same range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapa de Suportes")

# Clear the contents of AF2:AI71 (keeps number formatting / styles intact).
$ws.Range("AF2:AI71").ClearContents()

# Leave the selection where the user last worked, matching the cleared range.
[void]$ws.Range("AF2:AI71").Select()
